$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# Simple single-cell value edits
# ------------------------------------------------------------------
$ws.Cells.Item(20, 10).Value = 9     # J20: 20 -> 9
$ws.Cells.Item(30, 8).Value = 15     # H30: 10 -> 15
$ws.Cells.Item(30, 10).Value = 0     # J30: 10 -> 0

# ------------------------------------------------------------------
# Row 78/79 restructure:
#   old row 78 (PO 16063532 / "13.08.2021" / PCS / "10.09.2021") keeps its
#   B,C,D,F,G,H,J cells, but its A,E,I,K cells get overwritten with the
#   values that used to live in row 79 (PO 151442906 / 45350 / .. / 45348),
#   with the unit switching from PCS to M; the now-redundant row 79 is
#   removed outright, shifting every following row up by one.
# ------------------------------------------------------------------

# Grab the formatting (number format / font / fill / border) of cells that
# already look the way the target row 78 cells should look, before the old
# row 79 disappears underneath us.
$ws.Cells.Item(22, 1).Copy()                 # A22 has the "PO number" style
$ws.Cells.Item(78, 1).PasteSpecial(-4122)    # -4122 = xlPasteFormats

$ws.Cells.Item(79, 5).Copy()                 # E79 has the date style
$ws.Cells.Item(78, 5).PasteSpecial(-4122)

$ws.Cells.Item(79, 11).Copy()                # K79 has the date style
$ws.Cells.Item(78, 11).PasteSpecial(-4122)

# Now write the actual values into row 78.
$ws.Cells.Item(78, 1).Value = 151442906
$ws.Cells.Item(78, 5).Value = 45350
$ws.Cells.Item(78, 9).Value = "M"
$ws.Cells.Item(78, 11).Value = 45348

# Finally, delete the old row 79 entirely; rows 80/81 shift up to become 79/80.
$ws.Rows.Item(79).Delete()

# ------------------------------------------------------------------
# View state (matches the author's last on-screen position/selection)
# ------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 73
$ws.Range("J79").Select()
